# Add a new "2022-Q4" sheet right after "总计", shifting "2022-Q3" and
# "2022-Q2" one position to the right, and update the "总计" summary sheet
# with the new quarter's totals.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Update the "总计" (totals) sheet: 2022-Q4 becomes the new first data
#    row, and the previous two rows (2022-Q3 / 2022-Q2) shift down by one.
# ---------------------------------------------------------------------
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 13
$total.Cells.Item(2, 4).Value = 0.64

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2022-Q3"
$total.Cells.Item(3, 3).Value = 2
$total.Cells.Item(3, 4).Value = 0.01

$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = "2022-Q2"
$total.Cells.Item(4, 3).Value = 3
$total.Cells.Item(4, 4).Value = 0.27

# Match the existing label-column style (bold, centered, thin border)
$labelCol = $total.Range("A2:A4")
$labelCol.Font.Bold = $true
$labelCol.HorizontalAlignment = -4108
$labelCol.VerticalAlignment = -4160
$labelCol.Borders.LineStyle = 1
$labelCol.Borders.Weight = 2

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q4" worksheet right after "总计".
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

# Reproduce the page margins used by the other quarterly sheets (inches
# expressed in points: 1 in = 72 pt).
$q4.PageSetup.LeftMargin   = 54
$q4.PageSetup.RightMargin  = 54
$q4.PageSetup.TopMargin    = 72
$q4.PageSetup.BottomMargin = 72
$q4.PageSetup.HeaderMargin = 36
$q4.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------
# 3. Header row.
# ---------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $q4.Cells.Item(1, $col).Value = $headers[$col - 2]
}

$headerRange = $q4.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# ---------------------------------------------------------------------
# 4. Fund rows. Columns B-G are stored as text (even the numeric-looking
#    ones) except G13/G14 which are genuine zeros, matching the source
#    data export; A and H are real numbers throughout.
# ---------------------------------------------------------------------
$rows = @(
    @("012528", "广发鑫睿一年持有期混合A",       "3.82", "95.35", "5.67", "0.2166", 6),
    @("012529", "广发鑫睿一年持有期混合C",       "1.87", "95.35", "5.67", "0.1060", 6),
    @("000264", "博时内需增长混合A",             "2.31", "78.90", "4.55", "0.1051", 3),
    @("050012", "博时策略混合",                   "1.92", "79.11", "4.06", "0.0780", 6),
    @("166801", "浙商聚潮新思维混合A",           "1.46", "77.34", "4.49", "0.0656", 7),
    @("011351", "金鹰年年邮益一年持有期混合A",   "3.04", "39.17", "0.79", "0.0240", 6),
    @("530016", "建信恒稳价值混合",               "0.54", "52.32", "3.71", "0.0200", 1),
    @("001613", "长城久祥灵活配置混合A",         "0.24", "94.17", "5.71", "0.0137", 2),
    @("004677", "博时战略新兴产业混合",           "0.34", "86.55", "3.94", "0.0134", 4),
    @("011352", "金鹰年年邮益一年持有期混合C",   "0.23", "39.17", "0.79", "0.0018", 6),
    @("011982", "博时内需增长混合C",             "0.01", "78.90", "4.55", "0.0005", 3),
    @("014085", "浙商聚潮新思维混合C",           "0.00", "77.34", "4.49", $null,     7),
    @("017462", "长城久祥灵活配置混合C",         "0.00", "94.17", "5.71", $null,     2)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $fields = $rows[$i]

    $q4.Cells.Item($r, 1).Value = $i

    $q4.Cells.Item($r, 2).NumberFormat = "@"
    $q4.Cells.Item($r, 2).Value = $fields[0]

    $q4.Cells.Item($r, 3).Value = $fields[1]

    $q4.Cells.Item($r, 4).NumberFormat = "@"
    $q4.Cells.Item($r, 4).Value = $fields[2]

    $q4.Cells.Item($r, 5).NumberFormat = "@"
    $q4.Cells.Item($r, 5).Value = $fields[3]

    $q4.Cells.Item($r, 6).NumberFormat = "@"
    $q4.Cells.Item($r, 6).Value = $fields[4]

    if ($fields[5] -eq $null) {
        $q4.Cells.Item($r, 7).Value = 0
    } else {
        $q4.Cells.Item($r, 7).NumberFormat = "@"
        $q4.Cells.Item($r, 7).Value = $fields[5]
    }

    $q4.Cells.Item($r, 8).Value = $fields[6]
}

$labelColQ4 = $q4.Range("A2:A14")
$labelColQ4.Font.Bold = $true
$labelColQ4.HorizontalAlignment = -4108
$labelColQ4.VerticalAlignment = -4160
$labelColQ4.Borders.LineStyle = 1
$labelColQ4.Borders.Weight = 2

# ---------------------------------------------------------------------
# 5. Restore "2022-Q2" as the active/selected sheet, matching the
#    original workbook state (adding a sheet makes the new one active).
#    Re-resolve the sheet by name now, after the insert shifted indices.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Activate()
$q2.Range("A1").Select()
